$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "1"
$ws.Range("B2").Value = 50
$ws.Range("C2").Value = "ALL"

# Delete rows 3 and 4 entirely (shifting cells up)
$ws.Range("A3:C4").Delete()
